# Auto-generated Excel COM-interop script to apply scheduled market-price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Seraph Profits" sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Range("H9").Value = 218.75
$ws.Range("I9").Value = 260
$ws.Range("K9").Value = 260
$ws.Range("M9").Value = -91

# Row 17
$ws.Range("H17").Value = 1310.8823
$ws.Range("J17").Value = 1326.375
$ws.Range("L17").Value = 3979.125
$ws.Range("N17").Value = -4315.125

# Row 33
$ws.Range("H33").Value = 129
$ws.Range("I33").Value = 67.166664
$ws.Range("K33").Value = 67.166664
$ws.Range("M33").Value = 161.833336

# Row 98
$ws.Range("H98").Value = 1467.2142
$ws.Range("I98").Value = 1656.8
$ws.Range("K98").Value = 1656.8
$ws.Range("M98").Value = -158.8

# Row 122
$ws.Range("H122").Value = 1467.2142
$ws.Range("I122").Value = 1656.8
$ws.Range("K122").Value = 4970.4
$ws.Range("M122").Value = -2520.4

# Row 137
$ws.Range("H137").Value = 1540.4375
$ws.Range("I137").Value = 1387.3334
$ws.Range("K137").Value = 4162.0002
$ws.Range("M137").Value = -1612.0002


# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 500
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -288

# Row 97
$ws.Range("H97").Value = 1089.5
$ws.Range("I97").Value = 1073.7142
$ws.Range("K97").Value = 1073.7142
$ws.Range("M97").Value = -577.7141999999999

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# Row 136
$ws.Range("H136").Value = 500
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050


# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")

# Row 22
$ws.Range("H22").Value = 490.83334
$ws.Range("I22").Value = 390
$ws.Range("K22").Value = 390
$ws.Range("M22").Value = -217


# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$ws.Range("H7").Value = 97.07692
$ws.Range("I7").Value = 41.555557
$ws.Range("J7").Value = 222
$ws.Range("K7").Value = 41.555557
$ws.Range("L7").Value = 222
$ws.Range("M7").Value = 71.44444300000001
$ws.Range("N7").Value = -448

# Row 22
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -900

# Row 31
$ws.Range("H31").Value = 3490.2083
$ws.Range("I31").Value = 2457.8823
$ws.Range("K31").Value = 2457.8823
$ws.Range("M31").Value = -2162.8823

# Row 34
$ws.Range("H34").Value = 3490.2083
$ws.Range("I34").Value = 2457.8823
$ws.Range("K34").Value = 2457.8823
$ws.Range("M34").Value = -2255.8823

# Row 58
$ws.Range("H58").Value = 1985.1613
$ws.Range("I58").Value = 1092.375
$ws.Range("K58").Value = 1092.375
$ws.Range("M58").Value = -889.375

# Row 105
$ws.Range("H105").Value = 6772.5454
$ws.Range("I105").Value = 999.5
$ws.Range("J105").Value = 10071.429
$ws.Range("K105").Value = 999.5
$ws.Range("L105").Value = 10071.429
$ws.Range("M105").Value = 747.5
$ws.Range("N105").Value = -13565.429

# Row 122
$ws.Range("H122").Value = 3644.5
$ws.Range("I122").Value = 3644.5
$ws.Range("K122").Value = 10933.5
$ws.Range("M122").Value = -8483.5

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# Row 134
$ws.Range("H134").Value = 2935.3572
$ws.Range("I134").Value = 2791.6667
$ws.Range("J134").Value = 3194
$ws.Range("K134").Value = 8375.000100000001
$ws.Range("L134").Value = 9582
$ws.Range("M134").Value = -5840.000100000001
$ws.Range("N134").Value = -14652

# Row 136
$ws.Range("H136").Value = 1985.1613
$ws.Range("I136").Value = 1092.375
$ws.Range("K136").Value = 3277.125
$ws.Range("M136").Value = -727.125


# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")

# Row 8
$ws.Range("H8").Value = 198.33333
$ws.Range("I8").Value = 198.33333
$ws.Range("K8").Value = 594.99999
$ws.Range("M8").Value = -455.99999


# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

# Row 43
$ws.Range("H43").Value = 9797.4
$ws.Range("I43").Value = 9508.5
$ws.Range("K43").Value = 9508.5
$ws.Range("M43").Value = -9357.5

# Row 80
$ws.Range("H80").Value = 11000
$ws.Range("J80").Value = 11000
$ws.Range("L80").Value = 11000
$ws.Range("N80").Value = -12996

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 83
$ws.Range("H83").Value = 11000
$ws.Range("J83").Value = 11000
$ws.Range("L83").Value = 55000
$ws.Range("N83").Value = -64984

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 122
$ws.Range("H122").Value = 47515.863
$ws.Range("I122").Value = 1576.7222
$ws.Range("J122").Value = 254242
$ws.Range("K122").Value = 4730.1666
$ws.Range("L122").Value = 762726
$ws.Range("M122").Value = -2280.1666
$ws.Range("N122").Value = -767626


# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 55
$ws.Range("H55").Value = 778.8182
$ws.Range("I55").Value = 542.5714
$ws.Range("K55").Value = 542.5714
$ws.Range("M55").Value = -369.5714

# Row 82
$ws.Range("I82").Value = 4518.5
$ws.Range("K82").Value = 4518.5
$ws.Range("M82").Value = -4157.5

# Row 85
$ws.Range("I85").Value = 4518.5
$ws.Range("K85").Value = 4518.5
$ws.Range("M85").Value = -3270.5

# Row 122
$ws.Range("H122").Value = 24998.5
$ws.Range("I122").Value = 24998.5
$ws.Range("K122").Value = 74995.5
$ws.Range("M122").Value = -72545.5


# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")

# Row 62
$ws.Range("H62").Value = 7136.4546
$ws.Range("I62").Value = 3750
$ws.Range("K62").Value = 3750
$ws.Range("M62").Value = -3126

# Row 65
$ws.Range("H65").Value = 7136.4546
$ws.Range("I65").Value = 3750
$ws.Range("K65").Value = 18750
$ws.Range("M65").Value = -15630

# Row 100
$ws.Range("H100").Value = 4540
$ws.Range("I100").Value = 4675
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 9350
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = -8809
$ws.Range("N100").Value = -9082

# Row 122
$ws.Range("H122").Value = 2157.2856
$ws.Range("J122").Value = 1956.1111
$ws.Range("L122").Value = 5868.3333
$ws.Range("N122").Value = -10768.3333

# Row 132
$ws.Range("H132").Value = 3923.9048
$ws.Range("I132").Value = 3213.8572
$ws.Range("K132").Value = 9641.571599999999
$ws.Range("M132").Value = -7111.571599999999

# Row 136
$ws.Range("H136").Value = 1843.3572
$ws.Range("I136").Value = 1030.7
$ws.Range("K136").Value = 3092.1
$ws.Range("M136").Value = -542.1000000000004

